# Apply scheduled market-data refresh to the per-job Leve profit sheets.
# Each sheet has currentAveragePrice/NQ/HQ, LevePrice NQ/HQ and the derived
# LeveProfit NQ/HQ columns (H-N) refreshed from the latest market snapshot.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 1539.7858
$ws.Range("I33").Value = 403.3
$ws.Range("K33").Value = 403.3
$ws.Range("M33").Value = -174.3
$ws.Range("H51").Value = 14456.333
$ws.Range("I51").Value = 27125
$ws.Range("J51").Value = 4321.4
$ws.Range("K51").Value = 27125
$ws.Range("L51").Value = 4321.4
$ws.Range("M51").Value = -26641
$ws.Range("N51").Value = -5289.4
$ws.Range("H64").Value = 32187.117
$ws.Range("I64").Value = 201820
$ws.Range("J64").Value = 2940.0688
$ws.Range("K64").Value = 201820
$ws.Range("L64").Value = 2940.0688
$ws.Range("M64").Value = -201572
$ws.Range("N64").Value = -3436.0688
$ws.Range("H67").Value = 32187.117
$ws.Range("I67").Value = 201820
$ws.Range("J67").Value = 2940.0688
$ws.Range("K67").Value = 201820
$ws.Range("L67").Value = 2940.0688
$ws.Range("M67").Value = -200962
$ws.Range("N67").Value = -4656.0688
$ws.Range("H74").Value = 6233.3335
$ws.Range("I74").Value = 5850
$ws.Range("J74").Value = 7000
$ws.Range("K74").Value = 5850
$ws.Range("L74").Value = 7000
$ws.Range("M74").Value = -4914
$ws.Range("N74").Value = -8872
$ws.Range("H77").Value = 6233.3335
$ws.Range("I77").Value = 5850
$ws.Range("J77").Value = 7000
$ws.Range("K77").Value = 29250
$ws.Range("L77").Value = 35000
$ws.Range("M77").Value = -24570
$ws.Range("N77").Value = -44360
$ws.Range("H112").Value = 1199.375
$ws.Range("J112").Value = 1285
$ws.Range("L112").Value = 3855
$ws.Range("N112").Value = -6071
$ws.Range("H113").Value = 127438.125
$ws.Range("J113").Value = 2625
$ws.Range("L113").Value = 2625
$ws.Range("N113").Value = -9133
$ws.Range("H132").Value = 5819346
$ws.Range("I132").Value = 6102875
$ws.Range("J132").Value = 7000
$ws.Range("K132").Value = 18308625
$ws.Range("L132").Value = 21000
$ws.Range("M132").Value = -18306095
$ws.Range("N132").Value = -26060

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 30367.018
$ws.Range("I32").Value = 5276.18
$ws.Range("K32").Value = 5276.18
$ws.Range("M32").Value = -4989.18
$ws.Range("H55").Value = 10300
$ws.Range("J55").Value = 10485.714
$ws.Range("L55").Value = 10485.714
$ws.Range("N55").Value = -11115.714
$ws.Range("H132").Value = 2416.2354
$ws.Range("I132").Value = 2007.6364
$ws.Range("K132").Value = 6022.9092
$ws.Range("M132").Value = -3492.9092

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 111941.3
$ws.Range("I86").Value = 158358
$ws.Range("J86").Value = 3635.6667
$ws.Range("K86").Value = 158358
$ws.Range("L86").Value = 3635.6667
$ws.Range("M86").Value = -157235
$ws.Range("N86").Value = -5881.6667
$ws.Range("H89").Value = 111941.3
$ws.Range("I89").Value = 158358
$ws.Range("J89").Value = 3635.6667
$ws.Range("K89").Value = 791790
$ws.Range("L89").Value = 18178.3335
$ws.Range("M89").Value = -786174
$ws.Range("N89").Value = -29410.3335

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 1500
$ws.Range("I4").Value = 1500
$ws.Range("K4").Value = 1500
$ws.Range("M4").Value = -1388
$ws.Range("H21").Value = 0
$ws.Range("J21").Value = 0
$ws.Range("L21").Value = 0
$ws.Range("N21").ClearContents()
$ws.Range("H31").Value = 14879.156
$ws.Range("I31").Value = 21951.625
$ws.Range("J31").Value = 3173
$ws.Range("K31").Value = 21951.625
$ws.Range("L31").Value = 3173
$ws.Range("M31").Value = -21656.625
$ws.Range("N31").Value = -3763
$ws.Range("H34").Value = 14879.156
$ws.Range("I34").Value = 21951.625
$ws.Range("J34").Value = 3173
$ws.Range("K34").Value = 21951.625
$ws.Range("L34").Value = 3173
$ws.Range("M34").Value = -21749.625
$ws.Range("N34").Value = -3577
$ws.Range("H58").Value = 9237.727999999999
$ws.Range("I58").Value = 1678.7222
$ws.Range("J58").Value = 18308.533
$ws.Range("K58").Value = 1678.7222
$ws.Range("L58").Value = 18308.533
$ws.Range("M58").Value = -1475.7222
$ws.Range("N58").Value = -18714.533
$ws.Range("H86").Value = 2391.3333
$ws.Range("I86").Value = 2000
$ws.Range("J86").Value = 2782.6667
$ws.Range("K86").Value = 2000
$ws.Range("L86").Value = 2782.6667
$ws.Range("M86").Value = -877
$ws.Range("N86").Value = -5028.6667
$ws.Range("H89").Value = 2391.3333
$ws.Range("I89").Value = 2000
$ws.Range("J89").Value = 2782.6667
$ws.Range("K89").Value = 10000
$ws.Range("L89").Value = 13913.3335
$ws.Range("M89").Value = -4384
$ws.Range("N89").Value = -25145.3335
$ws.Range("H122").Value = 4109.32
$ws.Range("I122").Value = 3594.6667
$ws.Range("J122").Value = 5432.7144
$ws.Range("K122").Value = 10784.0001
$ws.Range("L122").Value = 16298.1432
$ws.Range("M122").Value = -8334.000100000001
$ws.Range("N122").Value = -21198.1432
$ws.Range("H136").Value = 9237.727999999999
$ws.Range("I136").Value = 1678.7222
$ws.Range("J136").Value = 18308.533
$ws.Range("K136").Value = 5036.1666
$ws.Range("L136").Value = 54925.599
$ws.Range("M136").Value = -2486.1666
$ws.Range("N136").Value = -60025.599

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 2227.182
$ws.Range("J34").Value = 2688.7778
$ws.Range("L34").Value = 8066.3334
$ws.Range("N34").Value = -8234.3334

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 1002561
$ws.Range("I5").Value = 2500400
$ws.Range("J5").Value = 4001.6667
$ws.Range("K5").Value = 2500400
$ws.Range("L5").Value = 4001.6667
$ws.Range("M5").Value = -2500288
$ws.Range("N5").Value = -4225.6667
$ws.Range("H70").Value = 91764.87
$ws.Range("I70").Value = 171074.83
$ws.Range("J70").Value = 5244.909
$ws.Range("K70").Value = 171074.83
$ws.Range("L70").Value = 5244.909
$ws.Range("M70").Value = -170804.83
$ws.Range("N70").Value = -5784.909
$ws.Range("H73").Value = 91764.87
$ws.Range("I73").Value = 171074.83
$ws.Range("J73").Value = 5244.909
$ws.Range("K73").Value = 171074.83
$ws.Range("L73").Value = 5244.909
$ws.Range("M73").Value = -170138.83
$ws.Range("N73").Value = -7116.909
$ws.Range("H132").Value = 2768.8948
$ws.Range("I132").Value = 2171.1
$ws.Range("J132").Value = 3433.111
$ws.Range("K132").Value = 6513.299999999999
$ws.Range("L132").Value = 10299.333
$ws.Range("M132").Value = -3983.299999999999
$ws.Range("N132").Value = -15359.333
$ws.Range("H134").Value = 28176.375
$ws.Range("J134").Value = 28176.375
$ws.Range("L134").Value = 84529.125
$ws.Range("N134").Value = -89599.125

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 75651.82000000001
$ws.Range("I2").Value = 267166.8
$ws.Range("J2").Value = 5585.3657
$ws.Range("K2").Value = 267166.8
$ws.Range("L2").Value = 5585.3657
$ws.Range("M2").Value = -267054.8
$ws.Range("N2").Value = -5809.3657
$ws.Range("H7").Value = 2087.652
$ws.Range("I7").Value = 1690.3529
$ws.Range("K7").Value = 1690.3529
$ws.Range("M7").Value = -1578.3529
$ws.Range("H16").Value = 101586.9
$ws.Range("I16").Value = 143409.72
$ws.Range("J16").Value = 4000.3333
$ws.Range("K16").Value = 143409.72
$ws.Range("L16").Value = 4000.3333
$ws.Range("M16").Value = -143239.72
$ws.Range("N16").Value = -4340.3333
$ws.Range("H46").Value = 844070
$ws.Range("I46").Value = 250
$ws.Range("J46").Value = 1446798.6
$ws.Range("K46").Value = 250
$ws.Range("L46").Value = 1446798.6
$ws.Range("M46").Value = -62
$ws.Range("N46").Value = -1447174.6
$ws.Range("H69").Value = 40000
$ws.Range("J69").Value = 0
$ws.Range("L69").Value = 0
$ws.Range("N69").ClearContents()
$ws.Range("H72").Value = 40000
$ws.Range("J72").Value = 0
$ws.Range("L72").Value = 0
$ws.Range("N72").ClearContents()
$ws.Range("H105").Value = 42000
$ws.Range("J105").Value = 42000
$ws.Range("L105").Value = 42000
$ws.Range("N105").Value = -48988
$ws.Range("H112").Value = 70799
$ws.Range("J112").Value = 70799
$ws.Range("L112").Value = 70799
$ws.Range("N112").Value = -73753
$ws.Range("H122").Value = 3152.5
$ws.Range("I122").Value = 5000
$ws.Range("J122").Value = 1305
$ws.Range("K122").Value = 15000
$ws.Range("L122").Value = 3915
$ws.Range("M122").Value = -12550
$ws.Range("N122").Value = -8815
$ws.Range("H126").Value = 2087.652
$ws.Range("I126").Value = 1690.3529
$ws.Range("K126").Value = 5071.0587
$ws.Range("M126").Value = -2601.0587
$ws.Range("H132").Value = 3378.1177
$ws.Range("I132").Value = 5102.933
$ws.Range("J132").Value = 2016.421
$ws.Range("K132").Value = 15308.799
$ws.Range("L132").Value = 6049.263
$ws.Range("M132").Value = -12778.799
$ws.Range("N132").Value = -11109.263
$ws.Range("H133").Value = 53433.332
$ws.Range("I133").Value = 45000
$ws.Range("J133").Value = 55120
$ws.Range("K133").Value = 45000
$ws.Range("L133").Value = 55120
$ws.Range("M133").Value = -42470
$ws.Range("N133").Value = -60180
$ws.Range("H135").Value = 40197.43
$ws.Range("J135").Value = 40197.43
$ws.Range("L135").Value = 40197.43
$ws.Range("N135").Value = -50337.43

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 25493.75
$ws.Range("J2").Value = 28990
$ws.Range("L2").Value = 28990
$ws.Range("N2").Value = -29214
$ws.Range("H81").Value = 500994.75
$ws.Range("J81").Value = 334659.66
$ws.Range("L81").Value = 669319.3199999999
$ws.Range("N81").Value = -671441.3199999999
$ws.Range("H84").Value = 500994.75
$ws.Range("J84").Value = 334659.66
$ws.Range("L84").Value = 3346596.6
$ws.Range("N84").Value = -3357204.6
$ws.Range("H132").Value = 3136.52
$ws.Range("I132").Value = 3412.7058
$ws.Range("J132").Value = 2549.625
$ws.Range("K132").Value = 10238.1174
$ws.Range("L132").Value = 7648.875
$ws.Range("M132").Value = -7708.117400000001
$ws.Range("N132").Value = -12708.875
$ws.Range("H136").Value = 1471.8948
$ws.Range("I136").Value = 609.5789
$ws.Range("J136").Value = 2334.2104
$ws.Range("K136").Value = 1828.7367
$ws.Range("L136").Value = 7002.6312
$ws.Range("M136").Value = 721.2633000000001
$ws.Range("N136").Value = -12102.6312
